# Update view-count figures (column F) on the "展览" and "全部类型" sheets.
# These two sheets both list the same events (the latter is an aggregate of
# all categories), so each changed event needs to be updated in both places.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value  = 611
$wsExhibition.Range("F5").Value  = 544
$wsExhibition.Range("F7").Value  = 2731
$wsExhibition.Range("F9").Value  = 7623
$wsExhibition.Range("F13").Value = 289

# Sheet "全部类型" (all types) - same events, different row positions
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 611
$wsAll.Range("F5").Value  = 544
$wsAll.Range("F9").Value  = 2731
$wsAll.Range("F11").Value = 7623
$wsAll.Range("F17").Value = 289
